$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before AM (so old AM/AN shift to AN/AO)
$ws.Range("AM1").EntireColumn.Insert()

# Set header text for the new AM1 cell
$ws.Range("AM1").Value = "antibodies_path"
$ws.Range("AN1").Value = "contributors_path"
$ws.Range("AO1").Value = "data_path"

# Make sure every header cell carries the correct comment after the shift.
if ($ws.Range("AM1").Comment -ne $null) { $ws.Range("AM1").Comment.Delete() }
if ($ws.Range("AN1").Comment -ne $null) { $ws.Range("AN1").Comment.Delete() }
if ($ws.Range("AO1").Comment -ne $null) { $ws.Range("AO1").Comment.Delete() }

$ws.Range("AM1").AddComment("Relative path to file with antibody information for this dataset.")
$ws.Range("AN1").AddComment("Relative path to file with ORCID IDs for contributors for this dataset.")
$ws.Range("AO1").AddComment("Relative path to file or directory with instrument data. Downstream processing will depend on filename extension conventions.")
